$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill")

# Copy formatting from row 9 to row 10, then fill in the new SKILL008 "blind" row
$ws.Range("A9:T9").Copy() | Out-Null
$ws.Range("A10:T10").PasteSpecial(-4122) | Out-Null
$ws.Rows("10:10").RowHeight = $ws.Rows("9:9").RowHeight

$ws.Range("A10").Value = "SKILL008"
$ws.Range("B10").Value = "目つぶし"
$ws.Range("C10").Value = "ATK_BST"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = "ENEMY_ONE"
$ws.Range("G10").Value = "AT_BLIND"
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("L10").Value = "blind"
$ws.Range("M10").Value = 60
$ws.Range("T10").Value = "相手を暗闇状態にする"
